$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading & Delivery Terms")
$ws.Columns("Q:Q").Delete()
